$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("training_subject_summary")
$ws.Range("C25").Value = "Completed"
$ws.Range("C28").Value = "In progress"
$ws.Range("C160").Value = "Completed"
$ws.Range("S160").Value = "not match"
$ws.Range("T160").Value = "not match"
$ws.Range("N163").Value = "not match"
$ws.Range("O163").Value = "not match"
$ws.Range("P163").Value = "not match"
$ws.Range("Q163").Value = "not match"
$ws.Range("R163").Value = "not match"
$ws.Range("S163").Value = "not match"
$ws.Range("M165").Value = "not match"
$ws.Range("N165").Value = "not match"
$ws.Range("P167").Value = "not match"
$ws.Range("Q167").Value = "not match"
$ws.Range("N168").Value = "not match"
$ws.Range("O168").Value = "not match"
$ws.Range("P168").Value = "not match"
$ws.Range("M169").Value = "not match"
$ws.Range("N169").Value = "not match"
$ws.Range("O169").Value = "not match"
$ws.Range("P169").Value = "not match"
$ws.Range("J170").Value = "not match"
$ws.Range("K170").Value = "not match"
$ws.Range("L170").Value = "not match"
$ws.Range("M170").Value = "not match"
$ws.Range("J171").Value = "not match"
$ws.Range("K171").Value = "not match"
$ws.Range("L171").Value = "not match"
$ws.Range("M171").Value = "not match"

$ws = $wb.Worksheets.Item("validation_subject_summary")
$ws.Range("F17").Value = 10
$ws.Range("H17").Value = 0
$ws.Range("F18").Value = 29
$ws.Range("H18").Value = 0
$ws.Range("F19").Value = 24
$ws.Range("H19").Value = 0
$ws.Range("F20").Value = 4
$ws.Range("H20").Value = 0
$ws.Range("F21").Value = 8
$ws.Range("H21").Value = 0
$ws.Range("F22").Value = 4
$ws.Range("H22").Value = 0
$ws.Range("F23").Value = 44
$ws.Range("H23").Value = 0
$ws.Range("F27").Value = 26
$ws.Range("H27").Value = 0
$ws.Range("D28").Value = 18
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 16
$ws.Range("H28").Value = 2
$ws.Range("M28").Value = "30-06-2023"
$ws.Range("U28").Value = 5
$ws.Range("D29").Value = 30
$ws.Range("E29").Value = 4
$ws.Range("F29").Value = 21
$ws.Range("O29").Value = "05-07-2023"
$ws.Range("P29").Value = "15-07-2023"
$ws.Range("U29").Value = 7
$ws.Range("D30").Value = 16
$ws.Range("E30").Value = 0
$ws.Range("H30").Value = 4
$ws.Range("L30").Value = "30-06-2023"
$ws.Range("U30").Value = 4
$ws.Range("D31").Value = 18
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 16
$ws.Range("H31").Value = 2
$ws.Range("M31").Value = "30-06-2023"
$ws.Range("U31").Value = 5
$ws.Range("D32").Value = 20
$ws.Range("E32").Value = 0
$ws.Range("F32").Value = 12
$ws.Range("L32").Value = "30-06-2023"
$ws.Range("M32").Value = "14-07-2023"
$ws.Range("U32").Value = 5
$ws.Range("F33").Value = 20
$ws.Range("H33").Value = 0
$ws.Range("F34").Value = 12
$ws.Range("H34").Value = 0
$ws.Range("D35").Value = 28
$ws.Range("E35").Value = 0
$ws.Range("F35").Value = 20
$ws.Range("H35").Value = 8
$ws.Range("N35").Value = "05-07-2023"
$ws.Range("O35").Value = "12-07-2023"
$ws.Range("U35").Value = 7
$ws.Range("D36").Value = 24
$ws.Range("E36").Value = 0
$ws.Range("F36").Value = 16
$ws.Range("H36").Value = 8
$ws.Range("M36").Value = "05-07-2023"
$ws.Range("N36").Value = "12-07-2023"
$ws.Range("U36").Value = 6
$ws.Range("F37").Value = 6
$ws.Range("H37").Value = 0
$ws.Range("D38").Value = 12
$ws.Range("E38").Value = 4
$ws.Range("F38").Value = 4
$ws.Range("H38").Value = 8
$ws.Range("K38").Value = "10-07-2023"
$ws.Range("L38").Value = "17-07-2023"
$ws.Range("U38").Value = 3
$ws.Range("F39").Value = 4
$ws.Range("H39").Value = 0
$ws.Range("D40").Value = 12
$ws.Range("E40").Value = 0
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = 12
$ws.Range("I40").Value = "05-07-2023"
$ws.Range("J40").Value = "12-07-2023"
$ws.Range("K40").Value = "19-07-2023"
$ws.Range("U40").Value = 3
$ws.Range("F41").Value = 13
$ws.Range("H41").Value = 0
$ws.Range("F42").Value = 23
$ws.Range("H42").Value = 0
$ws.Range("F43").Value = 17
$ws.Range("H43").Value = 0
$ws.Range("F44").Value = 24
$ws.Range("H44").Value = 0
$ws.Range("F45").Value = 49
$ws.Range("H45").Value = 0
$ws.Range("F46").Value = 22
$ws.Range("H46").Value = 0
$ws.Range("F47").Value = 11
$ws.Range("H47").Value = 0
$ws.Range("F48").Value = 44
$ws.Range("H48").Value = 0
$ws.Range("D52").Value = 26
$ws.Range("E52").Value = 0
$ws.Range("F52").Value = 12
$ws.Range("G52").Value = 4
$ws.Range("H52").Value = 10
$ws.Range("K52").Value = "19-07-2023"
$ws.Range("L52").Value = "11-08-2023"
$ws.Range("U52").Value = 4
$ws.Range("D53").Value = 9
$ws.Range("E53").Value = 6
$ws.Range("F53").Value = 4
$ws.Range("H53").Value = 5
$ws.Range("K53").Value = "22-08-2023"
$ws.Range("U53").Value = 2
$ws.Range("D54").Value = 4
$ws.Range("E54").Value = 7
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 0
$ws.Range("H54").Value = 4
$ws.Range("J54").Value = "19-07-2023"
$ws.Range("U54").Value = 1
$ws.Range("D55").Value = 8
$ws.Range("E55").Value = 5
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = 8
$ws.Range("J55").Value = "19-07-2023"
$ws.Range("K55").Value = "04-08-2023"
$ws.Range("U55").Value = 2
$ws.Range("D56").Value = 2
$ws.Range("E56").Value = 2
$ws.Range("F56").Value = 0
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 2
$ws.Range("J56").Value = "29-08-2023"
$ws.Range("U56").Value = 1
$ws.Range("F65").Value = 11
$ws.Range("H65").Value = 0
$ws.Range("D68").Value = 9
$ws.Range("E68").Value = 0
$ws.Range("F68").Value = 6
$ws.Range("H68").Value = 3
$ws.Range("J68").Value = "11-08-2023"
$ws.Range("U68").Value = 2
